# Update the consent-form heading:
#   "FECHA DE REALIZACIÓN DEL PROTOCOLO:" -> "FECHA DE REALIZACIÓN DEL CONSENTIMIENTO:"
#
# "PROTOCOLO" occurs exactly once in the document, inside the bold
# "REALIZACIÓN DEL PROTOCOLO" run, so a precise, case-sensitive,
# whole-word Find/Replace targets only that text without touching the
# neighboring "FECHA DE " / ":" runs.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Text = "PROTOCOLO"
$find.Replacement.Text = "CONSENTIMIENTO"
$find.Forward = $true
$find.Wrap = 0
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $true
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false

$find.Execute(
    "PROTOCOLO",    # FindText
    $true,          # MatchCase
    $true,          # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    "CONSENTIMIENTO", # ReplaceWith
    2               # Replace (wdReplaceAll)
)
